# "Use `<formatversion>` as suffix for table headers" / "Adjust `xlsx`
# export to new header formatting"
#
# The sheet has a header row (row 1) whose column names used to be suffixed
# generically with "_old" (columns A-J, the "before" half of the diff) and
# "_new" (columns L-U, the "after" half of the diff) - column K is just the
# literal "diff" marker column and is left alone. This edit renames those
# headers to use the concrete format-version identifiers instead
# ("_old" -> "_FV2410", "_new" -> "_FV2504"), wraps the whole A1:U70 range
# in a proper Excel Table so the new headers show up as table column names
# too, and freezes the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -------------------------------------
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$leftNames = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $leftCols.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $leftNames[$i]
}

# Column K ("diff") is unchanged.

$rightCols  = @("L","M","N","O","P","Q","R","S","T","U")
$rightNames = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)
for ($i = 0; $i -lt $rightCols.Length; $i++) {
    $ws.Range($rightCols[$i] + "1").Value = $rightNames[$i]
}

# --- 2. Convert A1:U70 into an Excel Table (ListObject) ---------------------
$tblRange = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $tblRange, $false, 1, "Table1")

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
